$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.168.60"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.916.53"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.7946"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.22%  "
$ws.Range("D6").Value = "'243.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.3190"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("D9").Value = "'26.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'0.06975"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "'0.08015"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'0.7551"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").Value = "1.911.67"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'5.240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'93.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "30.172.48"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'14.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'249.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.76%  "
$ws.Range("D20").Value = "'0.000007847"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "2.151.90"
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("D23").Value = "'0.9995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'6.995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'169.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "'9.355"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'0.1401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.53%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'2.059"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'1.393"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("D31").Value = "'1.529"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'4.377"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'4.133"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "'0.05422"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").Value = "'1.273"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").Value = "'0.01935"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").Value = "'2.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'6.195"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "'0.4473"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").Value = "'72.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").Value = "'1.914"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("D44").Value = "'0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'0.8350"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D46").Value = "'7.634"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").Value = "'9.854"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "'100.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "2.061.04"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'965.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.03%  "
$ws.Range("D51").Value = "'36.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
